$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Medium traffic density")

# Row 3: new NMAC-related data
$ws.Range("C3").Value = 196
$ws.Range("F3").Value = 76
$ws.Range("G3").Value = 76
$ws.Range("H3").Value = 36
$ws.Range("I3").Value = 6

# Row 4: new NMAC-related data
$ws.Range("C4").Value = 201
$ws.Range("F4").Value = 68
$ws.Range("G4").Value = 68
$ws.Range("H4").Value = 37
$ws.Range("I4").Value = 10

# Corrected NMAC counts for rows 9, 11, 12
$ws.Range("I9").Value = 13
$ws.Range("I11").Value = 14
$ws.Range("I12").Value = 10

# Row 17: new NMAC-related data
$ws.Range("C17").Value = 274
$ws.Range("F17").Value = 85
$ws.Range("G17").Value = 84
$ws.Range("H17").Value = 47
$ws.Range("I17").Value = 16

# Corrected NMAC counts for rows 25, 26
$ws.Range("I25").Value = 21
$ws.Range("I26").Value = 24

# Move the active selection to J26 on this sheet
$ws.Range("J26").Select()
